# Update creating order test cases and add test case for getting summary report
#
# Appends 19 new email-address test rows (for order creation / summary
# report test cases) to column A of Sheet1, directly below the existing
# 59 data rows (the sheet already has a header in A1 and data through
# A60).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newEmails = @(
    "ikk58050@zbock.com",
    "sns62539@zbock.com",
    "ool45473@omeie.com",
    "csz39978@nezid.com",
    "obo64382@omeie.com",
    "ssv86936@zbock.com",
    "abg80663@zbock.com",
    "uzm01189@zslsz.com",
    "xif30473@omeie.com",
    "klz21931@nezid.com",
    "gza24518@omeie.com",
    "tzl06814@nezid.com",
    "qky16640@zslsz.com",
    "hko75455@omeie.com",
    "spq38435@zbock.com",
    "icy48459@zbock.com",
    "dcn67384@zbock.com",
    "ywd16426@omeie.com",
    "awa67260@nezid.com"
)

$startRow = 61
for ($i = 0; $i -lt $newEmails.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newEmails[$i]
}
